# Word COM-interop script implementing the commit:
# "modification de quelques fautes orthographe + rajout de quelques notes a la fin"
#
# Summary of changes to word/document.xml:
#  1. Remove the old _GoBack bookmark that sits right after the threading
#     diagram picture (empty paragraph).
#  2. Split the big "Bilan personnel" paragraph into five separate
#     paragraphs, rewording/expanding the text about the weak points of
#     the group project and fixing "commiter" -> "commite".
#  3. Re-insert a _GoBack bookmark inside the new final paragraph, between
#     "autre" and "s" of "sur les autres".

$d = $word.ActiveDocument

# --- 1. remove the stale _GoBack bookmark near the picture -----------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. split paragraph 1 / 2 ----------------------------------------------
$d.Content.Find.Execute(
    "chacun. Le projet en lui",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "chacun. ^pLe projet en lui", 2) | Out-Null

# --- 3. split paragraph 2 / 3 ----------------------------------------------
$d.Content.Find.Execute(
    "(comme pour ce TP). Le point faible",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(comme pour ce TP). ^pLe point faible", 2) | Out-Null

# --- 4. reword end of paragraph 3, split paragraph 3 / 4 -------------------
$d.Content.Find.Execute(
    "partie réseau du TP qui était compliqué. Un deuxième",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "partie réseau du TP. Nous n’avons en effet pas eu le temps ni les compétences pour approfondir le projet sur ce domaine.^pUn deuxième", 2) | Out-Null

# split "partie " / "réseau du TP...domaine." into two separate runs
# (matches the two-run shape of the target, with no visible formatting change)
$r = $d.Content
$r.Find.Execute(
    "réseau du TP. Nous n’avons en effet pas eu le temps ni les compétences pour approfondir le projet sur ce domaine.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Bold = 1
$r.Bold = 0

# --- 5. reword/split paragraph 4 / 5 ----------------------------------------
$d.Content.Find.Execute(
    "de l’équipe. C’est-à-dire que, déjà, un membre de l’équipe est parti de l’école sans prévenir et sans «",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "de l’équipe.^pEn effet, un membre a quitté l’école sans avoir «", 2) | Out-Null

# --- 6. spelling fix: commiter -> commite -----------------------------------
$d.Content.Find.Execute(
    "commiter",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "commité", 2) | Out-Null

# --- 7. reword the tail of the final paragraph ------------------------------
$d.Content.Find.Execute(
    "» sa partie. Puis, d’autres, ne faisant pas leurs parties du travail qui leurs était demander.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "» sa partie du travail et sans qu’on puisse le contacter, de plus certains n’ont pas fait leur part de travail, ce qui a considérablement alourdit la charge de travail sur les autres et ce qui a rendu plus difficile la livraison du projet dans les délais.", 2) | Out-Null

# --- 8. re-insert the _GoBack bookmark between "autre" and "s" -------------
$r = $d.Content
$r.Find.Execute("sur les autres et ce qui a rendu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchor = $r.Duplicate
$anchor.Start = $r.End - ("s et ce qui a rendu".Length)
$anchor.End = $anchor.Start
$d.Bookmarks.Add("_GoBack", $anchor) | Out-Null
